# Update attendance / price figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5414
$ws1.Range("F7").Value = 118
$ws1.Range("F8").Value = 5394
$ws1.Range("F9").Value = 631
$ws1.Range("F12").Value = 13
$ws1.Range("G12").Value = 29.9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 5414
$ws4.Range("F8").Value = 118
$ws4.Range("F9").Value = 5394
$ws4.Range("F10").Value = 631
$ws4.Range("F13").Value = 13
$ws4.Range("G13").Value = 29.9
